# g1.9 - atualização na filtragem de ano e remoção da coluna categoria
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update rows 2-7 ("Variação em dez anos" section): new ordering, new values,
# and period changed from "2012 / 2022" to "2013 / 2022"
$dataDezAnos = @(
    @("Eletricidade e gás, água, esgoto, atividades de gestão de resíduos e descontaminação", 30.3618246924122),
    @("Atividades financeiras, de seguros e serviços relacionados", 27.4933798284287),
    @("Atividades imobiliárias", 23.40818553225802),
    @("Informação e comunicação", 9.265557954021503),
    @("Administração, defesa, educação e saúde públicas e seguridade social", 2.709036955545471),
    @("Agropecuária", -4.003418719932625)
)

$row = 2
foreach ($item in $dataDezAnos) {
    $ws.Cells.Item($row, 1).Value = $item[0]
    $ws.Cells.Item($row, 2).Value = $item[1]
    $ws.Cells.Item($row, 3).Value = "2013 / 2022"
    $row++
}

# Remove column D (Categoria) entirely - this shifts dimension from A1:D13 to A1:C13
$ws.Range("D1:D13").Delete()
